$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 634, pushing the existing rows 634-675 down to 635-676.
$ws.Rows.Item(634).Insert()

# Populate the newly inserted row 634 with the new weekly price record.
$ws.Range("A634").Value = 10
$ws.Range("B634").Value = "Vega Modelo de Temuco"
$ws.Range("C634").Value = "La Araucanía"
$ws.Range("D634").Value = 44746
$ws.Range("E634").Value = 9
$ws.Range("F634").Value = 100112003
$ws.Range("G634").Value = "Ajo"
$ws.Range("H634").Value = "Chino"
$ws.Range("I634").Value = "Primera"
$ws.Range("J634").Value = 450
$ws.Range("K634").Value = 20000
$ws.Range("L634").Value = 22000
$ws.Range("M634").Value = 21111
$ws.Range("N634").Value = "$/caja 10 kilos"
$ws.Range("O634").Value = "China"
$ws.Range("P634").Value = 2111
$ws.Range("Q634").Value = 10
$ws.Range("R634").Value = "Hortaliza"
